$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.587011218070984
$ws.Range("B1").Value = 4.395941734313965
$ws.Range("C1").Value = 3.040613412857056
$ws.Range("D1").Value = 1.197516560554504
$ws.Range("E1").Value = 0.912695050239563
